# Norway Eliteserien - base update (03-04-2024 22:09)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 183: odds tweak only ---
$ws.Cells.Item(183, 18).Value = 1.86   # R183
$ws.Cells.Item(183, 19).Value = 2.04   # S183

# --- Row 184: fixture replaced (Kristiansund BK vs KFUM) ---
$ws.Cells.Item(184, 2).Value = 7617324                 # B184 id
$ws.Cells.Item(184, 5).Value = 45389.39583333334        # E184 date
$ws.Cells.Item(184, 6).Value = "Kristiansund BK"        # F184 HomeTeam
$ws.Cells.Item(184, 7).Value = "KFUM"                   # G184 AwayTeam
$ws.Cells.Item(184, 11).Value = 2.375                   # K184
$ws.Cells.Item(184, 12).Value = 3.8                     # L184
$ws.Cells.Item(184, 13).Value = 2.6                     # M184
$ws.Cells.Item(184, 14).Value = 2.25                    # N184
$ws.Cells.Item(184, 15).Value = 3.8                     # O184
$ws.Cells.Item(184, 16).Value = 2.75                    # P184
$ws.Cells.Item(184, 17).Value = -0.25                   # Q184
$ws.Cells.Item(184, 18).Value = 2.02                    # R184
$ws.Cells.Item(184, 19).Value = 1.88                    # S184
$ws.Cells.Item(184, 20).Value = 2.75                    # T184
$ws.Cells.Item(184, 21).Value = 2                       # U184
$ws.Cells.Item(184, 22).Value = 1.85                    # V184

# --- Row 185: odds tweak only ---
$ws.Cells.Item(185, 18).Value = 1.84   # R185
$ws.Cells.Item(185, 19).Value = 2.06   # S185

# --- Row 186: fixture replaced (HamKam vs Molde) ---
$ws.Cells.Item(186, 2).Value = 7617325                 # B186 id
$ws.Cells.Item(186, 6).Value = "HamKam"                 # F186 HomeTeam
$ws.Cells.Item(186, 7).Value = "Molde"                  # G186 AwayTeam
$ws.Cells.Item(186, 11).Value = 4.8                     # K186
$ws.Cells.Item(186, 12).Value = 4.5                     # L186
$ws.Cells.Item(186, 13).Value = 1.571                   # M186
$ws.Cells.Item(186, 14).Value = 5.5                     # N186
$ws.Cells.Item(186, 15).Value = 4.75                    # O186
$ws.Cells.Item(186, 16).Value = 1.5                     # P186
$ws.Cells.Item(186, 17).Value = 1.25                    # Q186
$ws.Cells.Item(186, 18).Value = 1.85                    # R186
$ws.Cells.Item(186, 19).Value = 2.05                    # S186
$ws.Cells.Item(186, 20).Value = 3                       # T186
$ws.Cells.Item(186, 21).Value = 1.975                   # U186
$ws.Cells.Item(186, 22).Value = 1.875                   # V186

function Set-Row($r, $a, $b, $c, $d, $e, $f, $g, $k, $l, $m, $n, $o, $p, $q, $rr, $s, $t, $u, $v) {
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
    $ws.Cells.Item($r, 6).Value = $f
    $ws.Cells.Item($r, 7).Value = $g
    $ws.Cells.Item($r, 11).Value = $k
    $ws.Cells.Item($r, 12).Value = $l
    $ws.Cells.Item($r, 13).Value = $m
    $ws.Cells.Item($r, 14).Value = $n
    $ws.Cells.Item($r, 15).Value = $o
    $ws.Cells.Item($r, 16).Value = $p
    $ws.Cells.Item($r, 17).Value = $q
    $ws.Cells.Item($r, 18).Value = $rr
    $ws.Cells.Item($r, 19).Value = $s
    $ws.Cells.Item($r, 20).Value = $t
    $ws.Cells.Item($r, 21).Value = $u
    $ws.Cells.Item($r, 22).Value = $v
    $ws.Cells.Item($r, 23).Value = 0
    $ws.Cells.Item($r, 24).Value = 0
    $ws.Cells.Item($r, 25).Value = 0
    $ws.Cells.Item($r, 26).Value = 0
    $ws.Cells.Item($r, 27).Value = 0
}

Set-Row 187 185 7617328 "Norway Eliteserien" "Norway Eliteserien" 45389.5 "Sarpsborg" "Odd BK" 1.571 4.333 5 1.571 4.333 5 -1 1.95 1.95 3.25 1.925 1.925
Set-Row 188 186 7617329 "Norway Eliteserien" "Norway Eliteserien" 45389.5 "Stromsgodset" "Rosenborg" 2.3 3.75 2.8 2.6 3.75 2.45 0 2 1.9 2.75 1.875 1.975
Set-Row 189 187 7617330 "Norway Eliteserien" "Norway Eliteserien" 45389.59375 "SK Brann" "Fredrikstad" 1.4 5 6.5 1.4 5.5 6.5 -1.5 2 1.9 3.5 2.025 1.825

# Copy formatting (bold/border/center for col A, date format for col E) from row 186,
# reusing the existing style entries instead of inventing new ones.
$ws.Range("A186").Copy()
$ws.Range("A187:A189").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E186").Copy()
$ws.Range("E187:E189").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
